$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Market Size")
$ws.Copy($null, $ws)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Market Size till 2031"
$ws.Rows.Item(13).Delete()

Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) { Write-Host " -" $s.Name }

Write-Host "ws (orig) dim:" $ws.UsedRange.Address()
Write-Host "ws2 (copy) dim:" $ws2.UsedRange.Address()
